$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates
$ws.Range("D2").Value = "246.64"
$ws.Range("D4").Value = "5.280"
$ws.Range("D5").Value = "0.05800"
$ws.Range("D6").Value = "6.500"
$ws.Range("D7").Value = "3.134"
$ws.Range("D8").Value = "0.8164"
$ws.Range("D9").Value = "0.8578"
$ws.Range("D11").Value = "0.06937"
$ws.Range("D13").Value = "0.02878"
$ws.Range("D14").Value = "0.09405"
$ws.Range("D15").Value = "3.743"
$ws.Range("D16").Value = "0.001514"
$ws.Range("D17").Value = "0.04688"
$ws.Range("D18").Value = "0.0005956"
$ws.Range("D19").Value = "0.006268"
$ws.Range("D20").Value = "0.001235"
$ws.Range("D21").Value = "0.004619"
$ws.Range("D22").Value = "0.00006895"
$ws.Range("D24").Value = "2.142"
$ws.Range("D27").Value = "0.1358"
$ws.Range("D28").Value = "0.0002327"
$ws.Range("D40").Value = "0.03666"
$ws.Range("D41").Value = "0.006246"
$ws.Range("D43").Value = "0.003398"
$ws.Range("D44").Value = "0.007467"
$ws.Range("D45").Value = "0.00005257"
$ws.Range("D48").Value = "0.002250"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D50").Value = "0.0001999"

# Column E (Volume(1h)) updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
